# Auto-committed update: append a new "bringUpDate2" lookup definition
# to the DBS (Key/Search definitions) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# New row 6 - add the three values in the same order the source workbook's
# shared-string table picks them up (FunNm, then OrderBy, then ReadKey)
# so the resulting uniqueCount sequence matches: bringUpDate2,
# "BringUpDate DESC,CustId ASC", "BringUpDate <= , AND AcDate =".
$ws.Cells.Item(6, 1).Value = "bringUpDate2"
$ws.Cells.Item(6, 3).Value = "BringUpDate DESC,CustId ASC"
$ws.Cells.Item(6, 2).Value = "BringUpDate <= , AND AcDate ="

# Reflect the saved selection/scroll state: DBS sheet stays active with
# A7 selected (row below the newly added row).
$ws.Activate()
$ws.Range("A7").Select() | Out-Null
